$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($row, $col, $val) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

Set-TextCell 2 4 '27.058.06'
Set-TextCell 2 5 '  -1.03%  '

Set-TextCell 3 4 '1.821.99'
Set-TextCell 3 5 '  -0.85%  '

Set-TextCell 4 4 '1.013'
Set-TextCell 4 5 '  -0.28%  '

Set-TextCell 5 4 '311.59'
Set-TextCell 5 5 '  -1.07%  '

Set-TextCell 6 4 '1.011'
Set-TextCell 6 5 '  -0.23%  '

Set-TextCell 7 4 '0.4639'
Set-TextCell 7 5 '  -2.10%  '

Set-TextCell 8 4 '0.3626'
Set-TextCell 8 5 '  -1.95%  '

Set-TextCell 9 4 '0.07287'
Set-TextCell 9 5 '  -2.31%  '

Set-TextCell 10 4 '0.8660'
Set-TextCell 10 5 '  -2.08%  '

Set-TextCell 11 4 '19.84'
Set-TextCell 11 5 '  -3.14%  '

Set-TextCell 12 2 'WrappedEther'
Set-TextCell 12 3 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
Set-TextCell 12 4 '1.864.73'
Set-TextCell 12 5 '  -0.75%  '

Set-TextCell 13 2 'TRON'
Set-TextCell 13 3 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
Set-TextCell 13 4 '0.07608'
Set-TextCell 13 5 '  +2.95%  '

Set-TextCell 14 4 '92.70'
Set-TextCell 14 5 '  -0.52%  '

Set-TextCell 15 4 '5.331'
Set-TextCell 15 5 '  -2.69%  '

Set-TextCell 16 4 '6.457'
Set-TextCell 16 5 '  -1.81%  '

Set-TextCell 17 5 '  -0.43%  '

Set-TextCell 18 4 '0.000008620'
Set-TextCell 18 5 '  -2.46%  '

Set-TextCell 19 4 '1.011'
Set-TextCell 19 5 '  -0.30%  '

Set-TextCell 20 4 '27.376.19'
Set-TextCell 20 5 '  +0.04%  '

Set-TextCell 21 4 '14.46'
Set-TextCell 21 5 '  -2.55%  '

Set-TextCell 22 4 '5.181'
Set-TextCell 22 5 '  -3.20%  '

Set-TextCell 23 4 '10.55'
Set-TextCell 23 5 '  -1.41%  '

Set-TextCell 24 4 '2.096.28'
Set-TextCell 24 5 '  +0.40%  '

Set-TextCell 25 4 '151.42'
Set-TextCell 25 5 '  -0.56%  '

Set-TextCell 26 5 '  -2.56%  '

Set-TextCell 27 4 '18.22'
Set-TextCell 27 5 '  -2.30%  '

Set-TextCell 28 4 '2.098'
Set-TextCell 28 5 '  -3.29%  '

Set-TextCell 29 4 '115.91'
Set-TextCell 29 5 '  -1.71%  '

Set-TextCell 30 4 '5.070'
Set-TextCell 30 5 '  -3.68%  '

Set-TextCell 31 4 '0.08898'
Set-TextCell 31 5 '  -0.64%  '

Set-TextCell 32 4 '2.962'
Set-TextCell 32 5 '  +0.60%  '

Set-TextCell 33 4 '0.7291'
Set-TextCell 33 5 '  -3.99%  '

Set-TextCell 34 4 '4.435'
Set-TextCell 34 5 '  -2.68%  '

Set-TextCell 35 4 '1.137'
Set-TextCell 35 5 '  -3.40%  '

Set-TextCell 36 4 '1.011'
Set-TextCell 36 5 '  -0.30%  '

Set-TextCell 37 4 '2.535'
Set-TextCell 37 5 '  +6.64%  '

Set-TextCell 38 4 '1.072'
Set-TextCell 38 5 '  -3.19%  '

Set-TextCell 39 4 '0.05249'
Set-TextCell 39 5 '  -2.29%  '

Set-TextCell 40 2 'MXToken'
Set-TextCell 40 3 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
Set-TextCell 40 4 '2.945'
Set-TextCell 40 5 '  -1.95%  '

Set-TextCell 41 2 'VeChain'
Set-TextCell 41 3 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-TextCell 41 4 '0.01911'
Set-TextCell 41 5 '  -2.51%  '

Set-TextCell 42 4 '7.110'
Set-TextCell 42 5 '  -2.40%  '

Set-TextCell 43 4 '0.5199'
Set-TextCell 43 5 '  -2.83%  '

Set-TextCell 44 4 '0.1629'
Set-TextCell 44 5 '  -2.05%  '

Set-TextCell 45 4 '8.225'
Set-TextCell 45 5 '  -3.62%  '

Set-TextCell 46 4 '0.4850'
Set-TextCell 46 5 '  -2.45%  '

Set-TextCell 47 4 '1.011'
Set-TextCell 47 5 '  -0.37%  '

Set-TextCell 48 4 '103.19'
Set-TextCell 48 5 '  -1.98%  '

Set-TextCell 49 4 '10.06'
Set-TextCell 49 5 '  -4.58%  '

Set-TextCell 50 4 '1.638'
Set-TextCell 50 5 '  -2.38%  '

Set-TextCell 51 4 '0.06249'
Set-TextCell 51 5 '  -1.12%  '
